# Adding the changes we made on may 9th
# Insert 3 new accelerometer readings at the top (new rows 2-4, pushing the
# existing data down) and append 7 new readings at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new rows right after the header row (row 1) -----------------
$ws.Rows("2:4").Insert()
# Excel "insert" copies the formatting of the row above (the bold header);
# clear it so the new data rows look like plain, unformatted data rows.
$ws.Range("A2:C4").ClearFormats()

$topData = @(
    @(-0.8363723754882812, 5.859383583068848, 2.452773094177246),
    @(-0.9548721313476562, 5.95263671875, 2.935124397277832),
    @(-1.523673057556152, 5.953047752380371, 3.276031017303467)
)

$r = 2
foreach ($row in $topData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Append 7 new rows of readings at the bottom (rows 25-31) -------------
$bottomData = @(
    @(1.72843074798584, 5.596723556518555, -0.66290283203125),
    @(1.581844329833984, 5.306270599365234, -0.953785240650177),
    @(1.586828231811523, 5.404983997344971, -0.8601570129394531),
    @(1.426663398742676, 5.44196891784668, -0.7858069539070129),
    @(1.610628128051758, 5.431691646575928, -0.8632726669311523),
    @(1.379239082336426, 5.286327838897705, -0.7818757295608521),
    @(1.496992111206055, 5.291580200195312, -0.7745996713638306)
)

$r = 25
foreach ($row in $bottomData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
